$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "6.40")
# are preserved exactly as literal text instead of being parsed into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '63.274.86'
$ws.Range("E2").Value = '  +0.74%  '
$ws.Range("D3").Value = '3.026.30'
$ws.Range("E3").Value = '  -2.24%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '558.37'
$ws.Range("E5").Value = '  +0.39%  '
$ws.Range("D6").Value = '155.49'
$ws.Range("E6").Value = '  -3.65%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").Value = '0.560'
$ws.Range("E8").Value = '  -4.20%  '
$ws.Range("D9").Value = '3.029.85'
$ws.Range("E9").Value = '  -2.00%  '
$ws.Range("D10").Value = '0.113'
$ws.Range("E10").Value = '  -1.97%  '
$ws.Range("D11").Value = '6.40'
$ws.Range("E11").Value = '  -4.50%  '
$ws.Range("D12").Value = '0.367'
$ws.Range("E12").Value = '  -2.01%  '
$ws.Range("D13").Value = '3.556.61'
$ws.Range("E13").Value = '  -2.15%  '
$ws.Range("E14").Value = '  -3.04%  '
$ws.Range("D15").Value = '63.296.83'
$ws.Range("E15").Value = '  +0.41%  '
$ws.Range("D16").Value = '24.15'
$ws.Range("E16").Value = '  -1.02%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.0000151'
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.024.41'
$ws.Range("E18").Value = '  -2.40%  '
$ws.Range("D19").Value = '399.39'
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("D20").Value = '5.11'
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("D21").Value = '12.04'
$ws.Range("E21").Value = '  -2.30%  '
$ws.Range("D22").Value = '6.68'
$ws.Range("E22").Value = '  -4.20%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").Value = '65.42'
$ws.Range("E24").Value = '  -3.32%  '
$ws.Range("B25").Value = 'Kaspa'
$ws.Range("C25").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D25").Value = '0.190'
$ws.Range("E25").Value = '  -3.12%  '
$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D26").Value = '0.465'
$ws.Range("E26").Value = '  -1.86%  '
$ws.Range("D27").Value = '0.0₃0984'
$ws.Range("E27").Value = '  -1.35%  '
$ws.Range("D28").Value = '8.75'
$ws.Range("E28").Value = '  +2.00%  '
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").Value = '  -0.31%  '
$ws.Range("D31").Value = '1.76'
$ws.Range("E31").Value = '  -0.19%  '
$ws.Range("D32").Value = '20.38'
$ws.Range("E32").Value = '  -1.86%  '
$ws.Range("D33").Value = '162.90'
$ws.Range("E33").Value = '  +7.52%  '
$ws.Range("E34").Value = '  +3.15%  '
$ws.Range("D35").Value = '4.74'
$ws.Range("E35").Value = '  -0.53%  '
$ws.Range("D36").Value = '6.05'
$ws.Range("E36").Value = '  -1.43%  '
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("D38").Value = '2.546.07'
$ws.Range("E38").Value = '  -5.42%  '
$ws.Range("D39").Value = '1.60'
$ws.Range("E39").Value = '  -1.92%  '
$ws.Range("D40").Value = '22.93'
$ws.Range("E40").Value = '  -1.39%  '
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = '3.95'
$ws.Range("E41").Value = '  -0.97%  '
$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").Value = '37.78'
$ws.Range("E42").Value = '  -0.80%  '
$ws.Range("D43").Value = '0.668'
$ws.Range("E43").Value = '  -3.29%  '
$ws.Range("D44").Value = '0.0602'
$ws.Range("E44").Value = '  +0.51%  '
$ws.Range("D45").Value = '0.0251'
$ws.Range("E45").Value = '  -0.49%  '
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").Value = '5.11'
$ws.Range("E46").Value = '  -0.27%  '
$ws.Range("B47").Value = 'FirstDigitalUSD'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D47").Value = '0.997'
$ws.Range("E47").Value = '  -0.32%  '
$ws.Range("D48").Value = '20.39'
$ws.Range("E48").Value = '  -0.47%  '
$ws.Range("D49").Value = '271.08'
$ws.Range("E49").Value = '  -3.86%  '
$ws.Range("D50").Value = '10.48'
$ws.Range("E50").Value = '  +0.47%  '
$ws.Range("D51").Value = '0.0943'
$ws.Range("E51").Value = '  -2.28%  '

# Remove the temporary text-format styling from column D so the cells
# keep their original (unstyled) appearance, matching the source workbook.
$ws.Range("D2:D51").ClearFormats()
